$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Weekly Task Report")
$ws.Range("B7").Value = 11
